$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '61.738.39'
$ws.Range('E2').Value = '  +4.08%  '

# Row 3
$ws.Range('D3').Value = '3.085.13'
$ws.Range('E3').Value = '  +2.93%  '

# Row 4
$ws.Range('E4').Value = '  -0.10%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '577.72'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +2.31%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '142.29'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +2.24%  '

# Row 7
$ws.Range('E7').Value = '  -0.08%  '

# Row 8
$ws.Range('D8').Value = '3.078.01'
$ws.Range('E8').Value = '  +2.94%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.526'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +1.17%  '

# Row 10
$ws.Range('E10').Value = '  +4.24%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.48'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +8.71%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.467'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +2.11%  '

# Row 13
$ws.Range('E13').Value = '  +3.93%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '35.10'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +3.71%  '

# Row 15
$ws.Range('E15').Value = '  +0.12%  '

# Row 16
$ws.Range('D16').Value = '3.592.84'
$ws.Range('E16').Value = '  +2.91%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '7.27'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +3.73%  '

# Row 18
$ws.Range('D18').Value = '3.079.85'
$ws.Range('E18').Value = '  +2.80%  '

# Row 19
$ws.Range('D19').Value = '61.673.56'
$ws.Range('E19').Value = '  +3.98%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '449.84'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +4.60%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.95'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +2.34%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.731'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +2.52%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.45'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +3.67%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.59'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +0.73%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '82.10'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +1.63%  '

# Row 26
$ws.Range('E26').Value = '  +0.20%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.24'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +4.87%  '

# Row 28
$ws.Range('E28').Value = '  -0.15%  '

# Row 29
$ws.Range('E29').Value = '  +3.78%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.07'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +3.47%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.70'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +8.78%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '26.61'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +3.49%  '

# Row 33
$ws.Range('E33').Value = '  +9.33%  '

# Row 34
$ws.Range('E34').Value = '  +2.80%  '

# Row 35
$ws.Range('D35').Value = '0.0₃0795'
$ws.Range('E35').Value = '  +3.47%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.07'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +4.94%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.18'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +4.81%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '50.15'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +2.01%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.96'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +7.37%  '

# Row 40
$ws.Range('E40').Value = '  +1.90%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '430.25'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +6.46%  '

# Row 42
$ws.Range('E42').Value = '  +4.66%  '

# Row 43
$ws.Range('D43').Value = '2.789.53'
$ws.Range('E43').Value = '  +1.16%  '

# Row 44
$ws.Range('B44').Value = 'TheGraph'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.269'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +6.82%  '

# Row 45
$ws.Range('B45').Value = 'Kaspa'
$ws.Range('C45').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.108'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +0.14%  '

# Row 46
$ws.Range('B46').Value = 'Arweave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '35.63'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +10.22%  '

# Row 47
$ws.Range('B47').Value = 'Fetch.AI'
$ws.Range('C47').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.09'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +3.99%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '125.08'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +1.08%  '

# Row 49
$ws.Range('E49').Value = '  -0.02%  '

# Row 50
$ws.Range('E50').Value = '  +1.33%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '23.96'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +1.65%  '
